# BIS-1002: removed "Internal Assignment" column from export.
#
# The "Internal Assignment" header lives in column O (row 4), with the
# corresponding data values (FALSE) in column O of rows 5-8. Clearing the
# whole O4:O8 range removes the cell contents; since "Internal Assignment"
# then becomes an unused shared string it is dropped from the workbook's
# shared string table automatically on save, shifting every subsequent
# shared-string index down by one (exactly matching the target export).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O4:O8").ClearContents()

# Reflect the new selection left behind after removing the last column.
$ws.Range("O4:O8").Select()
